# ---------------------------------------------------------------------------
# Applies the "climate" sheet / planetary-boundaries restructuring commit to
# the Forrest_loss workbook.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forrest_loss")

# ---------------------------------------------------------------------------
# 1. Forrest_loss sheet edits
# ---------------------------------------------------------------------------

# Row 7: C7 now pulls from the wilderness table instead of the flat-rate calc
$ws.Range("C7").Formula = "=F15*B5"

# Row 9: label shifts down one slot in the (now shorter) shared-string table;
# the text itself is unchanged ("Has been suggested that is 3x the needed amount")
$ws.Range("A9").Value = "Has been suggested that is 3x the needed amount"

# Row 10 (new): relates the "Norway is paying" baseline to the low (primary forest)
# planetary-boundaries estimate converted to pounds
$ws.Range("C10").Formula = "=B8*C36/C8"

# Row 11: "Other wilderness:" - text unchanged, index shifts
$ws.Range("A11").Value = "Other wilderness:"

# Row 12 header row
$ws.Range("A12").Value = "Wilderness"
$ws.Range("B12").Value = "remaining"
$ws.Range("C12").Value = "lost 1990 - 2015"
$ws.Range("D12").Value = "protected"
$ws.Range("E12").Value = "taking $3 per km2 of remaining wilderness"
$ws.Range("F12").Value = "taking $10k per km2 of remaining wilderness"
$ws.Range("G12").Value = "taking £10k per km2 lost"
$ws.Range("H12").Value = "taking $950k per km2 lost"
$ws.Range("I12").Value = "assuming desert and tundra cost 1/10 to protect than rainforest"

# Rows 13-26: G and H columns get extra divisors (0.74 / 2); E13:E26 becomes a
# single shared-formula block (value-for-value identical to before)
$ws.Range("E13:E26").Formula = "=B13*3000000"
$ws.Range("G13").Formula = "=10000*1000000*C13/25/0.74"
$ws.Range("G14:G26").Formula = "=10000*1000000*C14/25/0.74"
$ws.Range("H13").Formula = "=950000*1000000*C13/25/2"
$ws.Range("H14:H26").Formula = "=950000*1000000*C14/25/2"

# Row 31 header: text unchanged ("low"/"mid"/"high"), shared-string indices shift
$ws.Range("B31").Value = "low"
$ws.Range("C31").Value = "mid"
$ws.Range("D31").Value = "high"

# Row 32: re-labelled "Primary forest loss" (was "Conclusion") with new formulas
$ws.Range("A32").Value = "Primary forest loss"
$ws.Range("B32").Formula = "=C6"
$ws.Range("C32").Formula = "=C8/2"
$ws.Range("D32").Formula = "=C8*2"

# Row 33 (new): Secondary forest loss
$ws.Range("A33").Value = "Secondary forest loss"
$ws.Range("B33").Formula = "=D6/2"
$ws.Range("C33").Formula = "=D6"
$ws.Range("D33").Formula = "=D6*2"

# Row 34 (new): Other land changes
$ws.Range("A34").Value = "Other land changes"
$ws.Range("B34").Formula = "=G28"
$ws.Range("C34").Formula = "=I28"
$ws.Range("D34").Formula = "=I28*4"

# Row 35 (new): section label
$ws.Range("A35").Value = "in Pounds"

# Rows 36-38 (new): same three categories converted to GBP (*0.73)
$ws.Range("A36").Value = "Primary forest loss"
$ws.Range("B36").Formula = "=B32 *0.73"
$ws.Range("C36:D36").Formula = "=C32 *0.73"

$ws.Range("A37").Value = "Secondary forest loss"
$ws.Range("B37:D37").Formula = "=B33 *0.73"

$ws.Range("A38").Value = "Other land changes"
$ws.Range("B38:D38").Formula = "=B34 *0.73"

# ---------------------------------------------------------------------------
# 2. New "climate" worksheet (placed after Forrest_loss)
# ---------------------------------------------------------------------------

$climate = $wb.Worksheets.Add($null, $ws)
$climate.Name = "climate"

$climate.Range("A1").Value = "Carbon estimates"
$climate.Range("B1").Value = "tCO2"
$climate.Range("D1").Value = "in Pounds"

$climate.Range("A2").Value = "WWF from consumption"
$climate.Range("B2").Formula = "=800/66.65"
$climate.Range("D2").Value = 33000000000

$climate.Range("A3").Value = "Woldwite production"
$climate.Range("B3").Formula = "=D2/7700000000"

$climate.Range("A4").Value = "Equal worldwide proportion"
$climate.Range("B4").Formula = "=D2*Forrest_loss!B5"

$climate.Range("A6").Value = "Per income worldwide proportion"
$climate.Range("B6").Value = "offsets"
$climate.Range("C6").Value = "cheap (7.4 pounds/tonne)"

$climate.Range("B7:B9").FormulaArray = "=B2:B4*7.4"
$climate.Range("C7:C9").FormulaArray = "=B2:B4*12.6"
$climate.Range("D7").Formula = "=B2*15*0.74"

$climate.Range("B8").Value = 31.714285714285715
$climate.Range("C8").Value = 54
$climate.Range("F8").Formula = "=15*0.74"

$climate.Range("B9").Value = 127.72032554932919
$climate.Range("C9").Value = 217.46974350291185

$climate.Columns.Item(1).ColumnWidth = 22.7109375
$climate.Columns.Item(4).ColumnWidth = 12

# ---------------------------------------------------------------------------
# 3. Final view state
# ---------------------------------------------------------------------------

[void]$climate.Range("F8").Select()
[void]$ws.Activate()
[void]$ws.Range("E35").Select()
